# Change n6 gateway IP from 10.0.2.15 to 10.0.2.2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 8 holds "global.n6network.gatewayIP" -> update both its columns (B8, E8)
$ws.Range("B8").Value = "10.0.2.2"
$ws.Range("E8").Value = "10.0.2.2"

# Move the active selection to the cell that was last edited
$ws.Range("E8").Select()
